# BOM update (inductore wirewound + RF in 1005)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill in the new BOM line (row 21): quantity, reference, description, unit price
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = "994-0402DC-27NXGRW "
$ws.Range("F21").Value = "RF inductor"
$ws.Range("G21").Value = 1.41

# Match the formatting used by the other data rows (D/F/G use the explicit
# black-font style already present on row 10/11, E stays with the default style)
$ws.Range("D21").Font.Color = $ws.Range("D10").Font.Color
$ws.Range("F21").Font.Color = $ws.Range("F10").Font.Color
$ws.Range("G21").Font.Color = $ws.Range("G11").Font.Color

# Recalculate so the shared formula in I21 (Prix total) reflects the new price
$excel.Calculate()

# Update the saved selection to match where the user left off editing
$ws.Range("G22").Select()
